$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$newText = @'
questions = [
    {
        "title": "Jamie is working as a storekeeper at your company on a monthly salary of $2,000. Along with recording movement of inventory, he is also responsible for conducting physical inventory counts at the end of each month and updating the records accordingly. At year-end, external auditors identified that the inventory in store was lower than the inventory shown in the accounts.Which of the following actions should you take?",
        "ques_type": 2,
        "options": [
            "Assign the monthly inventory count to Jamie\u2019s co-worker.",
            "Arrange for the monthly inventory count to be carried out by the external auditors and Jamie at the same time.",
            "Arrange a training session for the auditors to enhance their understanding of your company\u2019s business.",
            "Increase Jamie's salary by 25%."
        ],
        "score": "Assign the monthly inventory count to Jamie\u2019s co-worker."
    },
    {
        "title": "You are working as an operations manager at a zoo. Entry to the zoo is gained via two entrances, each with their own queue. However, due to a seasonal shortage of ticketing staff, visitors in both queues have to wait for up to 15 minutes before entering. As a result, many visitors are getting frustrated and you have received many complaints. Which of the following measures should you take?",
        "ques_type": 15,
        "options": [
            "Add 'estimated queue time from this point' signs to the queuing area.",
            "Allow visitors to switch between queues.",
            "Allow visitors who threaten to register formal complaints to enter first.",
            "Install video screens showing videos of different animals in the zoo.",
            "Allow visitors to enter without security checks."
        ],
        "score": [
            "Add 'estimated queue time from this point' signs to the queuing area.",
            "Install video screens showing videos of different animals in the zoo."
        ]
    },
    {
        "title": "Your company is facing production inefficiencies despite the fact that it uses the latest technology, promotes a healthy work environment, and has well-defined procedures. In an attempt to improve this situation, management is currently planning to run a recruitment drive.Which other option should you suggest?",
        "ques_type": 2,
        "options": [
            "\u201cLet's ensure we are fully compliant in terms of corporate governance.\u201d",
            "\u201cLet's ensure we are fully compliant in terms of labor laws.\u201d",
            "\u201cLet's introduce reward and bonus programs for the production staff.\u201d",
            "\u201cLet's increase our expectations of each staff member's workload.\u201d"
        ],
        "score": "\u201cLet's introduce reward and bonus programs for the production staff.\u201d"
    },
    {
        "title": "You are planning the procurement of various materials for next year. Your colleagues have information regarding material X in the table below. How many kilograms of material X should you plan to procure for next year?",
        "ques_type": 2,
        "options": [
            "686 kgs",
            "717 kgs",
            "914 kgs",
            "950 kgs"
        ],
        "score": "950 kgs"
    }
]
'@
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText
$ws.Rows("2:2").Delete()
$ws.Rows("1:1").EntireRow.AutoFit()
